$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "marks" column header in F1, matching the bold/centered header
# style already used by A1:E1 but with a left+right thin border instead of
# a full box border.
$ws.Range("F1").Value = "marks"
$ws.Range("F1").Font.Bold = $true
$ws.Range("F1").HorizontalAlignment = -4108
$ws.Range("F1").VerticalAlignment = -4160
$ws.Range("F1").Borders.Item(7).LineStyle = 1
$ws.Range("F1").Borders.Item(10).LineStyle = 1

# Student Alice Smith's marks value
$ws.Range("F2").Value = 1327

# Match the author's last on-screen selection after the edit
$ws.Range("F2").Select()
